# Spec_mul.xlsx edit script
# - Adds a new worksheet "Sheet5" at the end of the workbook with a small
#   data table (headers + 4 rows of measurement data).
# - Updates selections on a couple of the existing sheets.
# - Leaves the previously active sheet ("Agilent 34401A") selected/active.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Sheet5" worksheet after the last existing sheet.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)

$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "Sheet5"

# Header row (goes into sharedStrings.xml as new unique strings).
$ws5.Range("A1").Value = "x_"
$ws5.Range("B1").Value = "U_G_"
$ws5.Range("C1").Value = "U_FS_"
$ws5.Range("D1").Value = "U_"
$ws5.Range("E1").Value = "u_"
$ws5.Range("A1:E1").NumberFormat = "@"

# Data rows.
$ws5.Range("A2").Value = 1.56
$ws5.Range("B2").Value = 0.00078
$ws5.Range("C2").Value = 0.0002
$ws5.Range("D2").Value = 0.00098
$ws5.Range("E2").Value = 0.0006282051282051283

$ws5.Range("A3").Value = 78.9
$ws5.Range("B3").Value = 0.03945
$ws5.Range("C3").Value = 0.02
$ws5.Range("D3").Value = 0.05945
$ws5.Range("E3").Value = 0.0007534854245880862

$ws5.Range("A4").Value = 81.23
$ws5.Range("B4").Value = 0.040615
$ws5.Range("C4").Value = 0.02
$ws5.Range("D4").Value = 0.060615
$ws5.Range("E4").Value = 0.0007462144527883787

$ws5.Range("A5").Value = 100.67
$ws5.Range("B5").Value = 0.050335
$ws5.Range("C5").Value = 0.02
$ws5.Range("D5").Value = 0.070335
$ws5.Range("E5").Value = 0.0006986689182477403

# Column widths (closest values reachable through the ColumnWidth property,
# which internally snaps to whole-pixel increments).
$ws5.Columns.Item(1).ColumnWidth = 4.877604166666667
$ws5.Columns.Item(2).ColumnWidth = 7.877604166666667
$ws5.Columns.Item(3).ColumnWidth = 5.877604166666667
$ws5.Columns.Item(4).ColumnWidth = 7.877604166666667
$ws5.Columns.Item(5).ColumnWidth = 14.877604166666666

# ---------------------------------------------------------------------
# 2. Update the selection on "Keysight U1253B" and "Fluke 189".
# ---------------------------------------------------------------------
$wsKeysight = $wb.Worksheets.Item("Keysight U1253B")
[void]$wsKeysight.Activate()
[void]$wsKeysight.Range("E8").Select()

$wsFluke = $wb.Worksheets.Item("Fluke 189")
[void]$wsFluke.Activate()
[void]$wsFluke.Range("E2:E8").Select()

# ---------------------------------------------------------------------
# 3. Restore "Agilent 34401A" as the active sheet/tab (unchanged state).
# ---------------------------------------------------------------------
$wsAgilent = $wb.Worksheets.Item("Agilent 34401A")
[void]$wsAgilent.Activate()
[void]$wsAgilent.Range("H18").Select()
